$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")
$ws.Activate()
Write-Host "ActiveSheet:" $wb.ActiveSheet.Name
Write-Host "ActiveWindow TopLeftCell:" $excel.ActiveWindow.TopLeftCell.Address()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B22").Select()
Write-Host "Selection:" $excel.Selection.Address()
